# [Fonds de solidarite] Add 2020-07-17 data
# Update "nombre_aides" (column C) and "montant_total" (column D) figures
# for the rows impacted by the 2020-07-17 data refresh.
#
# Values are written with a leading apostrophe so Excel stores them as text
# (matching the workbook's existing convention of keeping every data cell,
# including numeric-looking ones, as a string) rather than converting them
# into native numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  C = "155";  D = "341600.00" },
    @{ Row = 3;  C = "837";  D = "2221538.52" },
    @{ Row = 4;  C = "343";  D = "1185844.92" },
    @{ Row = 21; C = "47";   D = "115500.00" },
    @{ Row = 22; C = "296";  D = "791242.04" },
    @{ Row = 24; C = "31";   D = "122500.00" },
    @{ Row = 25; C = "17";   D = "43500.00" },
    @{ Row = 44; C = "8";    D = "25500.00" },
    @{ Row = 45; C = "34";   D = "135903.07" },
    @{ Row = 46; C = "12";   D = "67880.00" },
    @{ Row = 47; C = "15";   D = "88197.00" },
    @{ Row = 48; C = "63";   D = "165636.00" },
    @{ Row = 49; C = "439";  D = "1195175.06" },
    @{ Row = 50; C = "194";  D = "650061.15" },
    @{ Row = 51; C = "55";   D = "257877.00" },
    @{ Row = 52; C = "15";   D = "60000.00" },
    @{ Row = 53; C = "11";   D = "28720.65" },
    @{ Row = 55; C = "2481"; D = "5816748.44" },
    @{ Row = 57; C = "432";  D = "1536075.00" },
    @{ Row = 59; C = "170";  D = "356000.00" },
    @{ Row = 71; C = "193";  D = "454326.09" },
    @{ Row = 72; C = "763";  D = "2093563.23" },
    @{ Row = 73; C = "286";  D = "973266.79" },
    @{ Row = 74; C = "89";   D = "348000.00" },
    @{ Row = 75; C = "20";   D = "105383.20" },
    @{ Row = 76; C = "20";   D = "41500.00" }
)

foreach ($u in $updates) {
    $ws.Range("C" + $u.Row).Value = "'" + $u.C
    $ws.Range("D" + $u.Row).Value = "'" + $u.D
}
